# feat: add new row for brands
#
# The sheet used to contain three identical rows of sample data
# (Wawi / Dwi / test). Replace that with a proper 2-row table: a header
# row (brand_name / founded_date / owner) and one data row describing
# the "Keebs.id" brand, owned by "Wawi", founded 2024-08-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the extra third row first so none of the old placeholder text
# ("Dwi", "test") lingers in the shared-string table.
$ws.Rows(3).Delete()

# Founded date: written/number-formatted before the header text so the
# auto-fit column width reflects only the date value, not the (longer)
# header string.
$ws.Range("B2").Value = 45517
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Columns("B").AutoFit()

# Header row.
$ws.Range("C1").Value = "owner"
$ws.Range("B1").Value = "founded_date"
$ws.Range("A1").Value = "brand_name"

# Remaining data cells: brand name, then owner (reuses the original
# "Wawi" string already in the sheet).
$ws.Range("A2").Value = "Keebs.id"
$ws.Range("C2").Value = "Wawi"

$ws.PageSetup.Orientation = 1

$ws.Range("C2").Select()
